$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master")

# --- Clean slate: remove existing hyperlinks and cell content on the Master sheet ---
$ws.Hyperlinks.Delete()
$ws.Cells.Clear()

# --- Column widths (best-effort; cosmetic/view-only, values per target layout) ---
$ws.Columns.Item(1).ColumnWidth = 2
$ws.Columns.Item(2).ColumnWidth = 16.43
$ws.Columns.Item(3).ColumnWidth = 19.71
$ws.Columns.Item(4).ColumnWidth = 9.14
$ws.Columns.Item(5).ColumnWidth = 2

# --- Row 2: title + "Updated:" label ---
$ws.Range("B2").Value = "SECTORS"
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").Font.Underline = $true
$ws.Range("D2").Value = "Updated:"

# --- Row 3 ---
$ws.Range("B3").Value = "Communications:"
$ws.Range("F3").Value = "Indexes:"
$ws.Range("L3").Value = "Goals:"
$ws.Range("M3").Value = "Focus on Med-device watchlist and work on models for all companies > `$10b"

# --- Row 4 ---
$ws.Range("C4").Value = "Telecom"
$ws.Range("G4").Value = "VOO"
$ws.Range("M4").Value = "Focus on Fin-Tech watchlist and work on models for all companies > `$10b"

# --- Row 5 ---
$ws.Range("C5").Value = "Trad Media"
$ws.Range("G5").Value = "QQQ"

# --- Row 6 ---
$ws.Range("C6").Value = "Interactive Media"
$ws.Range("G6").Value = "DJIA"

# --- Row 7: Discretionary header ---
$ws.Range("B7").Value = " Discretionary:"

# --- Rows 8-11: Discretionary items ---
$ws.Range("C8").Value = "E-Commerce"
$ws.Range("C9").Value = "Vehicles"
$ws.Range("C10").Value = "Clothing & Footwear"
$ws.Range("C11").Value = "Home Builders"

# --- Row 12: Staples header ---
$ws.Range("B12").Value = "Staples:"

# --- Rows 13-15: Staples items ---
$ws.Range("C13").Value = "Beverage Companies"
$ws.Range("C14").Value = "Tobacco"
$ws.Range("C15").Value = "Personal Products"

# --- Row 16: Energy header ---
$ws.Range("B16").Value = "Energy:"

# --- Row 17: Financials header ---
$ws.Range("B17").Value = "Financials:"

# --- Rows 18-21: Financials items ---
$ws.Range("C18").Value = "Fin-Tech"
$ws.Range("C19").Value = "BB & IB Fin."
$ws.Range("C20").Value = "PE"
$ws.Range("C21").Value = "Insurance"

# --- Row 22: Health Care header ---
$ws.Range("B22").Value = "Health Care:"

# --- Row 23: Health Care item ---
$ws.Range("C23").Value = "Med-Devices"

# --- Row 24: Industrial header ---
$ws.Range("B24").Value = "Industrial:"

# --- Row 25: Industrial item ---
$ws.Range("C25").Value = "Areo & Defense"

# --- Row 26: Materials header ---
$ws.Range("B26").Value = "Materials:"

# --- Row 27: Materials item ---
$ws.Range("C27").Value = "Chemicals"

# --- Row 28: Real Estate header ---
$ws.Range("B28").Value = "Real Estate:"

# --- Rows 29-30: Real Estate items ---
$ws.Range("C29").Value = "Home REITs"
$ws.Range("C30").Value = "Office REITs"

# --- Row 31: Technology header ---
$ws.Range("B31").Value = "Technology:"

# --- Rows 32-33: Technology items ---
$ws.Range("C32").Value = "SAAS"
$ws.Range("C33").Value = "Semis & Hardware"

# --- Row 34: Utilities header ---
$ws.Range("B34").Value = "Utilities:"

# --- Row 35: Utilities item ---
$ws.Range("C35").Value = "Nuclear"

# --- Hyperlinked cells: reapply the custom blue-underline "Hyperlink" look ---
$ws.Range("C13").Style = "Hyperlink"
$ws.Range("C13").Font.Color = 15773696
$ws.Range("C18").Style = "Hyperlink"
$ws.Range("C18").Font.Color = 15773696
$ws.Range("C23").Style = "Hyperlink"
$ws.Range("C23").Font.Color = 15773696

# --- Recreate hyperlinks (external links to related workbooks) ---
$ws.Hyperlinks.Add($ws.Range("C13"), "Consumer%20Staples\000%20Beverage%20Company%20Master%20List%20000.xlsx")
$ws.Hyperlinks.Add($ws.Range("C18"), "Financials\000%20Fin%20Tech%20Master%20List%20000.xlsx")
$ws.Hyperlinks.Add($ws.Range("C23"), "Health%20Care\000%20Med-Devices%20Master%20List%20000.xlsx", "", "", "Med Devices")

# --- Selection matches the saved view ---
$ws.Range("C18").Select()
